$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.539.12"
$ws.Range("E2").Value = "  +3.94%  "
$ws.Range("D3").Value = "1.740.54"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4811"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2681"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06244"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("D10").Value = "1.738.88"
$ws.Range("E10").Value = "  +4.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07127"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.71%  "
$ws.Range("E12").Value = "  +8.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6216"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.71%  "
$ws.Range("E14").Value = "  +4.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "26.542.81"
$ws.Range("E17").Value = "  +3.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006894"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.77%  "
$ws.Range("D21").Value = "1.961.16"
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.587"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.906"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "135.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("E26").Value = "  +3.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.811"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.423"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.012"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.742"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07892"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04590"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.003"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6380"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9319"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "112.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.001"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.434"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01517"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.98%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.743"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3922"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.986"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1200"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.22%  "
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.936"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.258"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3448"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.67%  "
